# Update countries & provincias Spain
# - Refresh the "last updated" timestamp (12:03 -> 13:03)
# - Catar / Israel swap places (and refresh their figures)
# - Emiratos Arabes Unidos / Japon swap places (and refresh their figures)
# - Malta / Tanzania swap places (and refresh their figures)
# - Refresh case counters for several other countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowStats($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Last updated timestamp ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 13:03"

# --- Country name swaps (rows keep their position, labels change) -------
$ws.Range("A31").Value = "Catar"
$ws.Range("A32").Value = "Israel"

$ws.Range("A34").Value = "Emiratos Arabes Unidos"
$ws.Range("A35").Value = "Japon"

$ws.Range("A117").Value = "Malta"
$ws.Range("A118").Value = "Tanzania"

# --- Updated statistics ---------------------------------------------------
Set-RowStats 9   166199 47   135100 24106 1949 0  6993   # Alemania
Set-RowStats 13  99970  1323 80475  13155 2685 63 6340   # Iran
Set-RowStats 31  17142  951  1924   15206 72   0  12     # now Catar
Set-RowStats 32  16268  22   10223  5808  89   2  237    # now Israel
Set-RowStats 34  15192  462  3153   11893 1    9  146    # now Emiratos Arabes Unidos
Set-RowStats 35  15078  0    4156   10386 328  0  536    # now Japon
Set-RowStats 37  13837  325  5454   7556  244  9  827    # Rumania
Set-RowStats 51  6849   24   5889   864   27   1  96     # Australia
Set-RowStats 56  5153   100  1799   3174  1    1  180    # Marruecos
Set-RowStats 60  4160   111  1264   2867  40   0  29     # Kazajistan
Set-RowStats 76  1946   20   911    956   4    1  79     # Bosnia y Herzegovina
Set-RowStats 85  1445   6    244    1103  17   1  98     # Eslovenia
Set-RowStats 97  820    17   570    219   7    0  31     # Albania
Set-RowStats 103 741    1    201    515   43   0  25     # Libano
Set-RowStats 117 482    2    403    74    1    1  5      # now Malta
Set-RowStats 118 480    0    167    297   7    0  16     # now Tanzania
